$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 157 is updated from 45170 to 45174
$ws.Range("C2:C157").Value = 45174
